$wb = $excel.ActiveWorkbook

# Insert a new blank column before column N ("Late"/"Outstanding" columns shift right)
# on the "Repayment Schedule" sheet.
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsSchedule.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab and select R4 on it
# (mirrors activeTab moving from the Transactions sheet to Repayment Schedule,
# and the new selected cell on that sheet).
$wsSchedule.Activate()
$wsSchedule.Range("R4").Select()
